# Fix the "Unidad" field label: the trailing punctuation after the
# label should be a colon (":"), matching the other field labels in the
# document (e.g. "REFERENCIA:", "N°:"), instead of a semicolon (";").
#
# "UNIDAD; +++unidad+++"  ->  "UNIDAD: +++unidad+++"

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "UNIDAD; +++unidad+++",  # FindText
    $true,                    # MatchCase
    $false,                   # MatchWholeWord
    $false,                   # MatchWildcards
    $false,                   # MatchSoundsLike
    $false,                   # MatchAllWordForms
    $true,                    # Forward
    1,                        # Wrap (wdFindContinue)
    $false,                   # Format
    "UNIDAD: +++unidad+++",   # ReplaceWith
    2                         # Replace (wdReplaceAll)
)

Write-Output "UNIDAD label fix applied: $found"
